$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 88 <= original row 89 (cols B, F:AC); A unchanged
$ws.Range("B88").Value = 6376945
$ws.Range("F88").Value = "Parnu JK Vaprus"
$ws.Range("G88").Value = "Harju JK Laagri"
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = "D"
$ws.Range("K88").Value = 1.615
$ws.Range("L88").Value = 4
$ws.Range("M88").Value = 4.5
$ws.Range("N88").Value = 1.85
$ws.Range("O88").Value = 3.8
$ws.Range("P88").Value = 3.5
$ws.Range("Q88").Value = -0.5
$ws.Range("R88").Value = 1.875
$ws.Range("S88").Value = 1.925
$ws.Range("T88").Value = 2.5
$ws.Range("U88").Value = 1.75
$ws.Range("V88").Value = 1.95
$ws.Range("W88").Value = -1
$ws.Range("X88").Value = 2.8
$ws.Range("Y88").Value = -1
$ws.Range("Z88").Value = -1
$ws.Range("AA88").Value = 0.925
$ws.Range("AB88").Value = -1
$ws.Range("AC88").Value = 0.95

# Row 89 <= original row 88 (cols B, F:AC); A unchanged
$ws.Range("B89").Value = 6376947
$ws.Range("F89").Value = "JK Tammeka Tartu"
$ws.Range("G89").Value = "JK Tallinna Kalev"
$ws.Range("H89").Value = 2
$ws.Range("I89").Value = 7
$ws.Range("J89").Value = "A"
$ws.Range("K89").Value = 3.6
$ws.Range("L89").Value = 3.4
$ws.Range("M89").Value = 1.909
$ws.Range("N89").Value = 2.4
$ws.Range("O89").Value = 3.6
$ws.Range("P89").Value = 2.45
$ws.Range("Q89").Value = 0
$ws.Range("R89").Value = 1.875
$ws.Range("S89").Value = 1.925
$ws.Range("T89").Value = 2.75
$ws.Range("U89").Value = 1.975
$ws.Range("V89").Value = 1.825
$ws.Range("W89").Value = -1
$ws.Range("X89").Value = -1
$ws.Range("Y89").Value = 1.45
$ws.Range("Z89").Value = -1
$ws.Range("AA89").Value = 0.925
$ws.Range("AB89").Value = 0.9750000000000001
$ws.Range("AC89").Value = -1

# Row 104 <= original row 106 (cols B, F:AC); A unchanged
$ws.Range("B104").Value = 6533597
$ws.Range("F104").Value = "FC Kuressaare"
$ws.Range("G104").Value = "Parnu JK Vaprus"
$ws.Range("H104").Value = 1
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = "H"
$ws.Range("K104").Value = 2.5
$ws.Range("L104").Value = 3.4
$ws.Range("M104").Value = 2.5
$ws.Range("N104").Value = 2.15
$ws.Range("O104").Value = 3.6
$ws.Range("P104").Value = 2.875
$ws.Range("Q104").Value = -0.25
$ws.Range("R104").Value = 1.95
$ws.Range("S104").Value = 1.85
$ws.Range("T104").Value = 2.75
$ws.Range("U104").Value = 1.95
$ws.Range("V104").Value = 1.85
$ws.Range("W104").Value = 1.15
$ws.Range("X104").Value = -1
$ws.Range("Y104").Value = -1
$ws.Range("Z104").Value = 0.95
$ws.Range("AA104").Value = -1
$ws.Range("AB104").Value = -1
$ws.Range("AC104").Value = 0.8500000000000001

# Row 105 <= original row 104 (cols B, F:AC); A unchanged
$ws.Range("B105").Value = 6537957
$ws.Range("F105").Value = "FC Flora Tallinn"
$ws.Range("G105").Value = "JK Nomme Kalju"
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = "D"
$ws.Range("K105").Value = 1.4
$ws.Range("L105").Value = 4
$ws.Range("M105").Value = 7.5
$ws.Range("N105").Value = 1.5
$ws.Range("O105").Value = 4.2
$ws.Range("P105").Value = 5
$ws.Range("Q105").Value = -1
$ws.Range("R105").Value = 1.85
$ws.Range("S105").Value = 1.95
$ws.Range("T105").Value = 2.75
$ws.Range("U105").Value = 1.85
$ws.Range("V105").Value = 1.95
$ws.Range("W105").Value = -1
$ws.Range("X105").Value = 3.2
$ws.Range("Y105").Value = -1
$ws.Range("Z105").Value = -1
$ws.Range("AA105").Value = 0.95
$ws.Range("AB105").Value = -1
$ws.Range("AC105").Value = 0.95

# Row 106 <= original row 105 (cols B, F:AC); A unchanged
$ws.Range("B106").Value = 6537869
$ws.Range("F106").Value = "JK Tallinna Kalev"
$ws.Range("G106").Value = "JK Trans Narva"
$ws.Range("H106").Value = 5
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = "H"
$ws.Range("K106").Value = 1.6
$ws.Range("L106").Value = 4
$ws.Range("M106").Value = 4.5
$ws.Range("N106").Value = 1.65
$ws.Range("O106").Value = 4
$ws.Range("P106").Value = 4.333
$ws.Range("Q106").Value = -0.75
$ws.Range("R106").Value = 1.8
$ws.Range("S106").Value = 2
$ws.Range("T106").Value = 2.75
$ws.Range("U106").Value = 1.9
$ws.Range("V106").Value = 1.9
$ws.Range("W106").Value = 0.6499999999999999
$ws.Range("X106").Value = -1
$ws.Range("Y106").Value = -1
$ws.Range("Z106").Value = 0.8
$ws.Range("AA106").Value = -1
$ws.Range("AB106").Value = 0.8999999999999999
$ws.Range("AC106").Value = -1

# Row 120 <= original row 121 (cols B, F:AC); A unchanged
$ws.Range("B120").Value = 7721087
$ws.Range("F120").Value = "Paide Linnameeskond"
$ws.Range("G120").Value = "FC Flora Tallinn"
$ws.Range("H120").Value = 2
$ws.Range("I120").Value = 1
$ws.Range("J120").Value = "H"
$ws.Range("K120").Value = 2.2
$ws.Range("L120").Value = 3.3
$ws.Range("M120").Value = 2.8
$ws.Range("N120").Value = 1.85
$ws.Range("O120").Value = 3.6
$ws.Range("P120").Value = 3.4
$ws.Range("Q120").Value = -0.5
$ws.Range("R120").Value = 1.9
$ws.Range("S120").Value = 1.9
$ws.Range("T120").Value = 2.5
$ws.Range("U120").Value = 1.95
$ws.Range("V120").Value = 1.85
$ws.Range("W120").Value = 0.8500000000000001
$ws.Range("X120").Value = -1
$ws.Range("Y120").Value = -1
$ws.Range("Z120").Value = 0.8999999999999999
$ws.Range("AA120").Value = -1
$ws.Range("AB120").Value = 0.95
$ws.Range("AC120").Value = -1

# Row 121 <= original row 120 (cols B, F:AC); A unchanged
$ws.Range("B121").Value = 7721007
$ws.Range("F121").Value = "JK Trans Narva"
$ws.Range("G121").Value = "JK Tammeka Tartu"
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 5
$ws.Range("J121").Value = "A"
$ws.Range("K121").Value = 2.25
$ws.Range("L121").Value = 3.3
$ws.Range("M121").Value = 2.75
$ws.Range("N121").Value = 2.1
$ws.Range("O121").Value = 3.25
$ws.Range("P121").Value = 3
$ws.Range("Q121").Value = -0.25
$ws.Range("R121").Value = 1.875
$ws.Range("S121").Value = 1.925
$ws.Range("T121").Value = 2.5
$ws.Range("U121").Value = 1.825
$ws.Range("V121").Value = 1.975
$ws.Range("W121").Value = -1
$ws.Range("X121").Value = -1
$ws.Range("Y121").Value = 2
$ws.Range("Z121").Value = -1
$ws.Range("AA121").Value = 0.925
$ws.Range("AB121").Value = 0.825
$ws.Range("AC121").Value = -1
